$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 5, 8, 1, 6, -4, -2, 54, 5),
    @(2, 6, 7, 1, 6, -5, -1, 65, 5),
    @(3, 6, 9, 5, 4, -1, -5, 21, 5),
    @(4, 5, 6, 2, 3, -3, -3, 43, 5),
    @(5, 8, 9, 6, 5, -2, -4, 32, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

$ws.Range("I1").Select()

$wb.Save()
